$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.876574238873843
$ws.Range("C2").Value = 0.3280797911055231
$ws.Range("D2").Value = 0.3963092732523847
$ws.Range("E2").Value = 0.08203988591942935
$ws.Range("G2").Value = 0.002805036148611686
$ws.Range("I2").Value = 7.90635596266435
$ws.Range("J2").Value = 0.01464460683253677
$ws.Range("K2").Value = 2.045061748202556
$ws.Range("L2").Value = 0.694616801155945
$ws.Range("M2").Value = 0.5766854997099458
$ws.Range("B3").Value = 1.894530729083613
$ws.Range("C3").Value = 0.3135922243317282
$ws.Range("D3").Value = 0.3962176827289454
$ws.Range("E3").Value = 0.08245112671865673
$ws.Range("G3").Value = 0.002811001998680823
$ws.Range("I3").Value = 7.582193094253796
$ws.Range("J3").Value = 0.01385338549802739
$ws.Range("K3").Value = 2.04375332438272
$ws.Range("L3").Value = 0.6907322584186488
$ws.Range("M3").Value = 0.5785158769421272
$ws.Range("B4").Value = 1.907165124553046
$ws.Range("C4").Value = 0.3047988964764556
$ws.Range("D4").Value = 0.3963213170383142
$ws.Range("E4").Value = 0.0827234049123442
$ws.Range("G4").Value = 0.002814851446213313
$ws.Range("I4").Value = 7.382211151095277
$ws.Range("J4").Value = 0.01336206671574303
$ws.Range("K4").Value = 2.044471758697341
$ws.Range("L4").Value = 0.6885697899885628
$ws.Range("M4").Value = 0.5799961747480111
$ws.Range("B5").Value = 1.912717883418111
$ws.Range("C5").Value = 0.3012409828577916
$ws.Range("D5").Value = 0.3964037725330058
$ws.Range("E5").Value = 0.08283934552504579
$ws.Range("G5").Value = 0.00281646717169081
$ws.Range("I5").Value = 7.300471390676165
$ws.Range("J5").Value = 0.01316044630803503
$ws.Range("K5").Value = 2.045146649716912
$ws.Range("L5").Value = 0.6877445451456623
$ws.Range("M5").Value = 0.5806890090205528
$ws.Range("B6").Value = 1.913664310653985
$ws.Range("C6").Value = 0.3006517237260198
$ws.Range("D6").Value = 0.3964198943274937
$ws.Range("E6").Value = 0.08285889882914788
$ws.Range("G6").Value = 0.002816738308189657
$ws.Range("I6").Value = 7.286883507603648
$ws.Range("J6").Value = 0.0131268820256274
$ws.Range("K6").Value = 2.045281778590862
$ws.Range("L6").Value = 0.6876108946265873
$ws.Range("M6").Value = 0.5808094645580972
$ws.Range("B7").Value = 1.907238375286028
$ws.Range("C7").Value = 0.3047508105098018
$ws.Range("D7").Value = 0.3963222661650718
$ws.Range("E7").Value = 0.08272494832707267
$ws.Range("G7").Value = 0.002814873045780243
$ws.Range("I7").Value = 7.381109784586641
$ws.Range("J7").Value = 0.01335935330297744
$ws.Range("K7").Value = 2.044479314042945
$ws.Range("L7").Value = 0.6885584338015747
$ws.Range("M7").Value = 0.580005155802251
$ws.Range("B8").Value = 1.88243139162347
$ws.Range("C8").Value = 0.3230631408600857
$ws.Range("D8").Value = 0.3962445237216343
$ws.Range("E8").Value = 0.08217758501594741
$ws.Range("G8").Value = 0.002807054593731703
$ws.Range("I8").Value = 7.794775762709577
$ws.Range("J8").Value = 0.01437292441426052
$ws.Range("K8").Value = 2.044294406250941
$ws.Range("L8").Value = 0.6932311888442371
$ws.Range("M8").Value = 0.577242615002902
$ws.Range("B9").Value = 1.846575571143177
$ws.Range("C9").Value = 0.359796389948599
$ws.Range("D9").Value = 0.3973600875428502
$ws.Range("E9").Value = 0.08126055915141617
$ws.Range("G9").Value = 0.002793193562724695
$ws.Range("I9").Value = 8.598874969698073
$ws.Range("J9").Value = 0.01631781170488367
$ws.Range("K9").Value = 2.056036328694233
$ws.Range("L9").Value = 0.7041628225763645
$ws.Range("M9").Value = 0.5746556432634478
$ws.Range("B10").Value = 1.828066321159554
$ws.Range("C10").Value = 0.3873066931409426
$ws.Range("D10").Value = 0.3989528110503926
$ws.Range("E10").Value = 0.08068139446813483
$ws.Range("G10").Value = 0.002783895335344273
$ws.Range("I10").Value = 9.185922867908459
$ws.Range("J10").Value = 0.01772212364421222
$ws.Range("K10").Value = 2.072090511703465
$ws.Range("L10").Value = 0.7132764854152072
$ws.Range("M10").Value = 0.574484786149938
$ws.Range("B11").Value = 1.821355650162474
$ws.Range("C11").Value = 0.39994002638295
$ws.Range("D11").Value = 0.399845370379893
$ws.Range("E11").Value = 0.08043829822409698
$ws.Range("G11").Value = 0.002779855192009762
$ws.Range("I11").Value = 9.452312473686163
$ws.Range("J11").Value = 0.01835596187957123
$ws.Range("K11").Value = 2.081018079423615
$ws.Range("L11").Value = 0.7176585908524942
$ws.Range("M11").Value = 0.5747837801506535
$ws.Range("B12").Value = 1.81906100170545
$ws.Range("C12").Value = 0.4047413501402843
$ws.Range("D12").Value = 0.4002075221288095
$ws.Range("E12").Value = 0.08034916082413979
$ws.Range("G12").Value = 0.002778352383521125
$ws.Range("I12").Value = 9.553102152071062
$ws.Range("J12").Value = 0.01859528612650863
$ws.Range("K12").Value = 2.084633177825481
$ws.Range("L12").Value = 0.7193520166387799
$ws.Range("M12").Value = 0.5749512567824766
$ws.Range("B13").Value = 1.819544218326371
$ws.Range("C13").Value = 0.4037065228227448
$ws.Range("D13").Value = 0.4001284518079586
$ws.Range("E13").Value = 0.08036822856068859
$ws.Range("G13").Value = 0.002778674837474917
$ws.Range("I13").Value = 9.531398986887439
$ws.Range("J13").Value = 0.01854377400910678
$ws.Range("K13").Value = 2.083844162141759
$ws.Range("L13").Value = 0.7189857933994119
$ws.Range("M13").Value = 0.5749127733475206
$ws.Range("B14").Value = 1.821161922715362
$ws.Range("C14").Value = 0.4003346847513285
$ws.Range("D14").Value = 0.3998746806479545
$ws.Range("E14").Value = 0.08043090642135908
$ws.Range("G14").Value = 0.002779731012694372
$ws.Range("I14").Value = 9.460606179374565
$ws.Range("J14").Value = 0.01837566506694444
$ws.Range("K14").Value = 2.081310792704841
$ws.Range("L14").Value = 0.7177972279774423
$ws.Range("M14").Value = 0.5747964707269446
$ws.Range("B15").Value = 1.822184942360167
$ws.Range("C15").Value = 0.3982716050078636
$ws.Range("D15").Value = 0.3997223846580908
$ws.Range("E15").Value = 0.08046967808607342
$ws.Range("G15").Value = 0.002780381475595852
$ws.Range("I15").Value = 9.417232568376562
$ws.Range("J15").Value = 0.01827260340684944
$ws.Range("K15").Value = 2.079789584656965
$ws.Range("L15").Value = 0.717073628827265
$ws.Range("M15").Value = 0.5747323000508473
$ws.Range("B16").Value = 1.828539337353078
$ws.Range("C16").Value = 0.3864834883773938
$ws.Range("D16").Value = 0.3988978609601617
$ws.Range("E16").Value = 0.08069769029575546
$ws.Range("G16").Value = 0.002784163170015878
$ws.Range("I16").Value = 9.168501122580778
$ws.Range("J16").Value = 0.01768060220178924
$ws.Range("K16").Value = 2.071539832227415
$ws.Range("L16").Value = 0.7129948622580287
$ws.Range("M16").Value = 0.5744728336303524
$ws.Range("B17").Value = 1.832875860914584
$ws.Range("C17").Value = 0.3792824726925517
$ws.Range("D17").Value = 0.3984350742254037
$ws.Range("E17").Value = 0.08084277735320278
$ws.Range("G17").Value = 0.002786531576130013
$ws.Range("I17").Value = 9.015749267367198
$ws.Range("J17").Value = 0.0173161654066547
$ws.Range("K17").Value = 2.066895514391689
$ws.Range("L17").Value = 0.71055321168609
$ws.Range("M17").Value = 0.5744102026476057
$ws.Range("B18").Value = 1.83553095377934
$ws.Range("C18").Value = 0.3751518000295562
$ws.Range("D18").Value = 0.398184706499066
$ws.Range("E18").Value = 0.08092814544421145
$ws.Range("G18").Value = 0.002787911683030797
$ws.Range("I18").Value = 8.92782764727761
$ws.Range("J18").Value = 0.01710608001464564
$ws.Range("K18").Value = 2.064377073958127
$ws.Range("L18").Value = 0.7091710733068055
$ws.Range("M18").Value = 0.5744096321227019
$ws.Range("B19").Value = 1.836457524655287
$ws.Range("C19").Value = 0.3737551338298601
$ws.Range("D19").Value = 0.3981026526350604
$ws.Range("E19").Value = 0.08095737934675462
$ws.Range("G19").Value = 0.002788382036588931
$ws.Range("I19").Value = 8.898047845926243
$ws.Range("J19").Value = 0.01703486718100322
$ws.Range("K19").Value = 2.063550600374754
$ws.Range("L19").Value = 0.708706922509549
$ws.Range("M19").Value = 0.5744155256461738
$ws.Range("B20").Value = 1.832397579035273
$ws.Range("C20").Value = 0.3800478759872021
$ws.Range("D20").Value = 0.3984827018988426
$ws.Range("E20").Value = 0.08082713417870924
$ws.Range("G20").Value = 0.002786277607799066
$ws.Range("I20").Value = 9.032016398455994
$ws.Range("J20").Value = 0.01735500891918207
$ws.Range("K20").Value = 2.067374084826383
$ws.Range("L20").Value = 0.7108108280690146
$ws.Range("M20").Value = 0.5744131999777338
$ws.Range("B21").Value = 1.82068006669553
$ws.Range("C21").Value = 0.4013246020847703
$ws.Range("D21").Value = 0.3999485637298079
$ws.Range("E21").Value = 0.08041241731203197
$ws.Range("G21").Value = 0.002779420053867455
$ws.Range("I21").Value = 9.481402013647426
$ws.Range("J21").Value = 0.01842506147314182
$ws.Range("K21").Value = 2.082048535788658
$ws.Range("L21").Value = 0.7181454148489763
$ws.Range("M21").Value = 0.5748291584970602
$ws.Range("B22").Value = 1.814459368113376
$ws.Range("C22").Value = 0.4153314622887478
$ws.Range("D22").Value = 0.4010474064023555
$ws.Range("E22").Value = 0.08015837839090523
$ws.Range("G22").Value = 0.002775096164427995
$ws.Range("I22").Value = 9.774602760401592
$ws.Range("J22").Value = 0.01912035268772883
$ws.Range("K22").Value = 2.093005938602857
$ws.Range("L22").Value = 0.7231372956580202
$ws.Range("M22").Value = 0.5754172884639175
$ws.Range("B23").Value = 1.817647698486184
$ws.Range("C23").Value = 0.4078463759520616
$ws.Range("D23").Value = 0.4004480487402162
$ws.Range("E23").Value = 0.08029241165055456
$ws.Range("G23").Value = 0.002777389511700369
$ws.Range("I23").Value = 9.618158684071886
$ws.Range("J23").Value = 0.01874962622829202
$ws.Range("K23").Value = 2.087032421397367
$ws.Range("L23").Value = 0.7204548741559904
$ws.Range("M23").Value = 0.5750744241969343
$ws.Range("B24").Value = 1.83261330583602
$ws.Range("C24").Value = 0.3797018082513546
$ws.Range("D24").Value = 0.3984611205429474
$ws.Range("E24").Value = 0.08083420036575006
$ws.Range("G24").Value = 0.002786392369553554
$ws.Range("I24").Value = 9.024662348556774
$ws.Range("J24").Value = 0.01733744953132188
$ws.Range("K24").Value = 2.06715725087102
$ws.Range("L24").Value = 0.7106942924293378
$ws.Range("M24").Value = 0.5744117345130206
$ws.Range("B25").Value = 1.854902523088271
$ws.Range("C25").Value = 0.3497688600402569
$ws.Range("D25").Value = 0.3969225461347321
$ws.Range("E25").Value = 0.08149197759612115
$ws.Range("G25").Value = 0.002796787029884529
$ws.Range("I25").Value = 8.382042844470334
$ws.Range("J25").Value = 0.01579609306561736
$ws.Range("K25").Value = 2.051559088565
$ws.Range("L25").Value = 0.7010157971964048
$ws.Range("M25").Value = 0.5750520515516513
